$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Fill in "Have" (column C) quantities to match "Qty" (column B) for these rows
$ws.Range("C17").Value = 4
$ws.Range("C18").Value = 2
$ws.Range("C19").Value = 2
$ws.Range("C20").Value = 2
$ws.Range("C21").Value = 2
$ws.Range("C22").Value = 2
$ws.Range("C24").Value = 1
$ws.Range("C25").Value = 1

# Row 18 also has a "Bought" (column D) note marker
$ws.Range("D18").Value = "*"

# Move the active selection to A2
$ws.Range("A2").Select()
